{"js": "// 1) Update the experience-years figure in the Professional Summary.\nconst summaryHits = context.document.body.search(\"21 years of experience\", { matchCase: true });\nsummaryHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < summaryHits.items.length; i++) {\n  summaryHits.items[i].insertText(\"15+ years of experience\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove the EDUCATION section entirely: the \"EDUCATION\" Heading2\n//    paragraph plus the two Heading3 degree paragraphs that follow it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \"EDUCATION\" && p.style === \"Heading 2\") {\n    toDelete.push(p);\n    if (i + 1 < paragraphs.items.length) toDelete.push(paragraphs.items[i + 1]);\n    if (i + 2 < paragraphs.items.length) toDelete.push(paragraphs.items[i + 2]);\n    break;\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the experience-years figure in the Professional Summary.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"21 years of experience\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"15+ years of experience\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Remove the EDUCATION section: the \"EDUCATION\" Heading2 paragraph plus\n#    the two Heading3 degree paragraphs that immediately follow it.\n$eduIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text.Trim() -eq \"EDUCATION\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n    $eduIndex = $i\n    break\n  }\n}\n\nif ($eduIndex -gt 0) {\n  $startPara = $d.Paragraphs.Item($eduIndex)\n  $endPara = $d.Paragraphs.Item($eduIndex + 2)\n  $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n  $range.Delete()\n}\n"}
